$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19 (Excel copies formatting from the row above,
# shifting the old rows 19-60 down to 20-61).
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly data point.
$ws.Cells.Item(19, 1).Value = 4
$ws.Cells.Item(19, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(19, 3).Value = "Los Lagos"
$ws.Cells.Item(19, 4).Value = 44519
$ws.Cells.Item(19, 5).Value = 10
$ws.Cells.Item(19, 6).Value = 100112026
$ws.Cells.Item(19, 7).Value = "Haba"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 120
$ws.Cells.Item(19, 11).Value = 10000
$ws.Cells.Item(19, 12).Value = 10000
$ws.Cells.Item(19, 13).Value = 10000
$ws.Cells.Item(19, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(19, 15).Value = "Región del Maule"
$ws.Cells.Item(19, 16).Value = 400
$ws.Cells.Item(19, 17).Value = 25
$ws.Cells.Item(19, 18).Value = "Hortaliza"
